$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "LiDAR corps" worksheet after the last existing sheet
#    (after "Options"). Worksheets.Add(Before, After) mirrors real Excel's
#    COM signature; passing Missing for Before and the last sheet for After
#    appends it at the end, exactly like the author manually inserting a new
#    tab.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lidar = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$lidar.Name = "LiDAR corps"

# Populate manufacturer name (col A) + link (col B) rows. The link (col B)
# is written before the name (col A) for each row so the shared-string
# table order matches the original authoring order (url then name).
$lidar.Range("B1").Value = "https://www.sick.com/be/en/detection-and-ranging-solutions/3d-lidar-sensors/c/g282752"
$lidar.Range("A1").Value = "SICK"

$lidar.Range("B2").Value = "https://leddartech.com/"
$lidar.Range("A2").Value = "Leddartech"

$lidar.Range("B3").Value = "https://quanergy.com/"
$lidar.Range("A3").Value = "Quanergy"

$lidar.Range("B4").Value = "https://www.robosense.ai/"
$lidar.Range("A4").Value = "Robosense"

$lidar.Range("B5").Value = "https://velodynelidar.com/"
$lidar.Range("A5").Value = "Velodyne"

$lidar.Range("B6").Value = "https://www.neuvition.com/"
$lidar.Range("A6").Value = "Neuvition"

$lidar.Range("B7").Value = "https://www.teledyneoptech.com/en/home/"
$lidar.Range("A7").Value = "Teledyne"

$lidar.Range("B8").Value = "https://www.zxlidars.com/"
$lidar.Range("A8").Value = "ZX LiDAR"

$lidar.Range("B9").Value = "https://www.aeye.ai/"
$lidar.Range("A9").Value = "Aeye"

$lidar.Range("A10").Value = "/+ System interators"

$lidar.Range("A11").Value = "https://www.continental-automotive.com/en-gl/Passenger-Cars/Autonomous-Mobility/Enablers/Lidars"

# Column A width (manufacturer names column)
$lidar.Columns.Item(1).ColumnWidth = 9.92

# ---------------------------------------------------------------------------
# 2. "Prototype V1.0": widen column B and move the selection
# ---------------------------------------------------------------------------
$proto = $wb.Worksheets.Item("Prototype V1.0")
$proto.Activate() | Out-Null
$proto.Columns.Item(2).ColumnWidth = 29.25
$proto.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. "Options": turn the camera-sensor link (C2) into a real hyperlink and
#    move the selection
# ---------------------------------------------------------------------------
$options = $wb.Worksheets.Item("Options")
$options.Activate() | Out-Null
$c2 = $options.Range("C2")
$linkText = $c2.Text
$options.Hyperlinks.Add($c2, $linkText) | Out-Null
$options.Range("D30").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Leave "LiDAR corps" as the active tab with its own selection, matching
#    the author's final view state.
# ---------------------------------------------------------------------------
$lidar.Activate() | Out-Null
$lidar.Range("P14").Select() | Out-Null
